# Daily attendance processing - 2025-10-11 16:21:03
# Normalizes the "Recorded By" (column G) entries so the "System" / "system"
# token that was recorded last is moved to the front of the comma separated
# list, matching the canonical "System, <user>" ordering used elsewhere in
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows in column G whose "Recorded By" value needs its first and last
# comma-separated entries swapped (i.e. the trailing System/system marker
# is promoted to the front of the list).
$rows = @(2, 3, 4, 5, 6, 10, 12, 13, 14, 15, 29, 30, 31, 32, 33, 37, 39, 40, 41, 42, 56, 57, 58, 59, 60, 64, 66, 67, 68, 69, 84, 85, 86, 87, 88, 89, 93, 95, 110, 111, 112, 113, 114, 115, 119, 121, 136, 137, 138, 139, 140, 141, 145, 147)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $current = [string]$cell.Value2

    $rawParts = $current.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    if ($parts.Count -gt 1) {
        $first = $parts[0]
        $last = $parts[$parts.Count - 1]
        $parts[0] = $last
        $parts[$parts.Count - 1] = $first
    }

    $cell.Value = [string]::Join(", ", $parts)
}
